$d = $word.ActiveDocument

$replacements = @(
    @("524÷9=58, 2", "266÷6=44, 2"),
    @("612÷8=76, 4", "402÷7=57, 3"),
    @("339÷2=169, 1", "822÷7=117, 3"),
    @("972÷6=162, 0", "719÷7=102, 5"),
    @("444÷8=55, 4", "507÷6=84, 3"),
    @("878÷5=175, 3", "676÷8=84, 4"),
    @("364÷8=45, 4", "335÷2=167, 1"),
    @("662÷2=331, 0", "958÷7=136, 6"),
    @("529÷8=66, 1", "133÷5=26, 3"),
    @("227÷7=32, 3", "278÷8=34, 6"),
    @("513÷9=57, 0", "920÷4=230, 0"),
    @("657÷3=219, 0", "637÷5=127, 2"),
    @("148÷5=29, 3", "293÷4=73, 1"),
    @("356÷5=71, 1", "616÷7=88, 0"),
    @("836÷7=119, 3", "193÷7=27, 4"),
    @("364÷6=60, 4", "384÷7=54, 6"),
    @("930÷8=116, 2", "804÷2=402, 0"),
    @("521÷5=104, 1", "125÷6=20, 5"),
    @("519÷9=57, 6", "917÷5=183, 2"),
    @("589÷4=147, 1", "757÷5=151, 2"),
    @("274÷3=91, 1", "976÷9=108, 4"),
    @("499÷3=166, 1", "851÷7=121, 4"),
    @("689÷9=76, 5", "608÷3=202, 2"),
    @("236÷2=118, 0", "202÷8=25, 2"),
    @("742÷2=371, 0", "920÷7=131, 3")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
